$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1434.5454
$ws.Range("I40").Value = 1333.3334
$ws.Range("J40").Value = 1472.5
$ws.Range("K40").Value = 1333.3334
$ws.Range("L40").Value = 1472.5
$ws.Range("M40").Value = -1158.3334
$ws.Range("N40").Value = -1822.5

$ws.Range("H64").Value = 2911
$ws.Range("I64").Value = 2942.7144
$ws.Range("K64").Value = 2942.7144
$ws.Range("M64").Value = -2694.7144

$ws.Range("H67").Value = 2911
$ws.Range("I67").Value = 2942.7144
$ws.Range("K67").Value = 2942.7144
$ws.Range("M67").Value = -2084.7144

$ws.Range("H76").Value = 3579.1428
$ws.Range("I76").Value = 3275
$ws.Range("J76").Value = 3984.6667
$ws.Range("K76").Value = 3275
$ws.Range("L76").Value = 3984.6667
$ws.Range("M76").Value = -2960
$ws.Range("N76").Value = -4614.6667

$ws.Range("H79").Value = 3579.1428
$ws.Range("I79").Value = 3275
$ws.Range("J79").Value = 3984.6667
$ws.Range("K79").Value = 3275
$ws.Range("L79").Value = 3984.6667
$ws.Range("M79").Value = -2183
$ws.Range("N79").Value = -6168.6667

$ws.Range("H116").Value = 632797.25
$ws.Range("I116").Value = 911340.4
$ws.Range("K116").Value = 911340.4
$ws.Range("M116").Value = -907898.4

$ws.Range("H129").Value = 800.8182
$ws.Range("J129").Value = 972.25
$ws.Range("L129").Value = 2916.75
$ws.Range("N129").Value = -12916.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 8660620
$ws.Range("I63").Value = 23087568
$ws.Range("J63").Value = 4451
$ws.Range("K63").Value = 23087568
$ws.Range("L63").Value = 4451
$ws.Range("M63").Value = -23086882
$ws.Range("N63").Value = -5823

$ws.Range("H66").Value = 8660620
$ws.Range("I66").Value = 23087568
$ws.Range("J66").Value = 4451
$ws.Range("K66").Value = 115437840
$ws.Range("L66").Value = 22255
$ws.Range("M66").Value = -115434408
$ws.Range("N66").Value = -29119

$ws.Range("H74").Value = 1872.3939
$ws.Range("I74").Value = 1338.2759
$ws.Range("K74").Value = 1338.2759
$ws.Range("M74").Value = -464.2759000000001

$ws.Range("H77").Value = 1872.3939
$ws.Range("I77").Value = 1338.2759
$ws.Range("K77").Value = 6691.379500000001
$ws.Range("M77").Value = -2323.379500000001

$ws.Range("H80").Value = 38241.5
$ws.Range("J80").Value = 38241.5
$ws.Range("L80").Value = 38241.5
$ws.Range("N80").Value = -40237.5

$ws.Range("H83").Value = 38241.5
$ws.Range("J83").Value = 38241.5
$ws.Range("L83").Value = 114724.5
$ws.Range("N83").Value = -124708.5

$ws.Range("H88").Value = 9525712
$ws.Range("I88").Value = 16668266
$ws.Range("J88").Value = 2308
$ws.Range("K88").Value = 16668266
$ws.Range("L88").Value = 2308
$ws.Range("M88").Value = -16667860
$ws.Range("N88").Value = -3120

$ws.Range("H91").Value = 9525712
$ws.Range("I91").Value = 16668266
$ws.Range("J91").Value = 2308
$ws.Range("K91").Value = 16668266
$ws.Range("L91").Value = 2308
$ws.Range("M91").Value = -16666862
$ws.Range("N91").Value = -5116

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("N89").ClearContents()

$ws.Range("H105").Value = 2494.8333
$ws.Range("I105").Value = 2293.2222
$ws.Range("J105").Value = 3099.6667
$ws.Range("K105").Value = 2293.2222
$ws.Range("L105").Value = 3099.6667
$ws.Range("M105").Value = -546.2222000000002
$ws.Range("N105").Value = -6593.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 55561130
$ws.Range("I62").Value = 166671600
$ws.Range("J62").Value = 5892
$ws.Range("K62").Value = 166671600
$ws.Range("L62").Value = 5892
$ws.Range("M62").Value = -166670976
$ws.Range("N62").Value = -7140

$ws.Range("H65").Value = 55561130
$ws.Range("I65").Value = 166671600
$ws.Range("J65").Value = 5892
$ws.Range("K65").Value = 833358000
$ws.Range("L65").Value = 29460
$ws.Range("M65").Value = -833354880
$ws.Range("N65").Value = -35700

$ws.Range("H130").Value = 42746.668
$ws.Range("J130").Value = 42746.668
$ws.Range("L130").Value = 42746.668
$ws.Range("N130").Value = -52786.668

$ws.Range("H132").Value = 4069.5173
$ws.Range("I132").Value = 3717.5293
$ws.Range("J132").Value = 4568.1665
$ws.Range("K132").Value = 11152.5879
$ws.Range("L132").Value = 13704.4995
$ws.Range("M132").Value = -8622.5879
$ws.Range("N132").Value = -18764.4995

$ws.Range("H134").Value = 4126.475
$ws.Range("I134").Value = 4401.1333
$ws.Range("J134").Value = 3302.5
$ws.Range("K134").Value = 13203.3999
$ws.Range("L134").Value = 9907.5
$ws.Range("M134").Value = -10668.3999
$ws.Range("N134").Value = -14977.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 51062.8
$ws.Range("J107").Value = 127123.875
$ws.Range("L107").Value = 381371.625
$ws.Range("N107").Value = -385211.625

$ws.Range("H113").Value = 543.6129
$ws.Range("I113").Value = 559
$ws.Range("K113").Value = 1677
$ws.Range("M113").Value = 493

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6499.407
$ws.Range("I70").Value = 5865.5264
$ws.Range("J70").Value = 8004.875
$ws.Range("K70").Value = 5865.5264
$ws.Range("L70").Value = 8004.875
$ws.Range("M70").Value = -5595.5264
$ws.Range("N70").Value = -8544.875

$ws.Range("H73").Value = 6499.407
$ws.Range("I73").Value = 5865.5264
$ws.Range("J73").Value = 8004.875
$ws.Range("K73").Value = 5865.5264
$ws.Range("L73").Value = 8004.875
$ws.Range("M73").Value = -4929.5264
$ws.Range("N73").Value = -9876.875

$ws.Range("H80").Value = 35716784
$ws.Range("I80").Value = 125001250
$ws.Range("K80").Value = 125001250
$ws.Range("M80").Value = -125000252

$ws.Range("H83").Value = 35716784
$ws.Range("I83").Value = 125001250
$ws.Range("K83").Value = 625006250
$ws.Range("M83").Value = -625001258

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3960.6553
$ws.Range("I132").Value = 1608.4849
$ws.Range("J132").Value = 7065.52
$ws.Range("K132").Value = 4825.4547
$ws.Range("L132").Value = 21196.56
$ws.Range("M132").Value = -2295.4547
$ws.Range("N132").Value = -26256.56

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H127").Value = 39830
$ws.Range("J127").Value = 39830
$ws.Range("L127").Value = 39830
$ws.Range("N127").Value = -49750

$ws.Range("H132").Value = 6292803
$ws.Range("I132").Value = 4066.862
$ws.Range("J132").Value = 13891693
$ws.Range("K132").Value = 12200.586
$ws.Range("L132").Value = 41675079
$ws.Range("M132").Value = -9670.585999999999
$ws.Range("N132").Value = -41680139
